# Auto-generated edit script applying the Golem_Profits.xlsx diff
# across all affected sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 90000
$ws.Range("J3").Value = 90000
$ws.Range("L3").Value = 90000
$ws.Range("N3").Value = -90228
$ws.Range("H7").Value = 2276
$ws.Range("J7").Value = 1750
$ws.Range("L7").Value = 1750
$ws.Range("N7").Value = -1974
$ws.Range("H14").Value = 2276
$ws.Range("J14").Value = 1750
$ws.Range("L14").Value = 1750
$ws.Range("N14").Value = -2132
$ws.Range("H40").Value = 3339.6
$ws.Range("I40").Value = 2900
$ws.Range("J40").Value = 3449.5
$ws.Range("K40").Value = 2900
$ws.Range("L40").Value = 3449.5
$ws.Range("M40").Value = -2725
$ws.Range("N40").Value = -3799.5
$ws.Range("H58").Value = 99
$ws.Range("I58").Value = 99
$ws.Range("K58").Value = 297
$ws.Range("M58").Value = -147
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H102").Value = 90000
$ws.Range("J102").Value = 90000
$ws.Range("L102").Value = 90000
$ws.Range("N102").Value = -96490
$ws.Range("H106").Value = 1000000000
$ws.Range("I106").Value = 1000000000
$ws.Range("K106").Value = 1000000000
$ws.Range("M106").Value = -999999369
$ws.Range("H131").Value = 4388.75
$ws.Range("I131").Value = 499.5
$ws.Range("J131").Value = 8278
$ws.Range("K131").Value = 1498.5
$ws.Range("L131").Value = 24834
$ws.Range("M131").Value = 3541.5
$ws.Range("N131").Value = -34914
$ws.Range("H132").Value = 597.0476
$ws.Range("I132").Value = 614.4
$ws.Range("K132").Value = 1843.2
$ws.Range("M132").Value = 686.8000000000002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4953.0586
$ws.Range("I32").Value = 4953.0586
$ws.Range("K32").Value = 4953.0586
$ws.Range("M32").Value = -4666.0586
$ws.Range("H45").Value = 2265.1667
$ws.Range("I45").Value = 1918.2
$ws.Range("J45").Value = 4000
$ws.Range("K45").Value = 1918.2
$ws.Range("L45").Value = 4000
$ws.Range("M45").Value = -1541.2
$ws.Range("N45").Value = -4754

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 474.75
$ws.Range("I24").Value = 466.33334
$ws.Range("J24").Value = 500
$ws.Range("K24").Value = 466.33334
$ws.Range("L24").Value = 500
$ws.Range("M24").Value = -231.33334
$ws.Range("N24").Value = -970
$ws.Range("H102").Value = 10459.167
$ws.Range("I102").Value = 10459.167
$ws.Range("K102").Value = 10459.167
$ws.Range("M102").Value = -7214.166999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 500008
$ws.Range("H27").Value = 500008
$ws.Range("H39").Value = 2550.3333
$ws.Range("I39").Value = 2550.3333
$ws.Range("K39").Value = 2550.3333
$ws.Range("M39").Value = -2159.3333
$ws.Range("H49").Value = 2550.3333
$ws.Range("I49").Value = 2550.3333
$ws.Range("K49").Value = 2550.3333
$ws.Range("M49").Value = -2368.3333
$ws.Range("H58").Value = 946.7778
$ws.Range("I58").Value = 953.1429000000001
$ws.Range("J58").Value = 924.5
$ws.Range("K58").Value = 953.1429000000001
$ws.Range("L58").Value = 924.5
$ws.Range("M58").Value = -750.1429000000001
$ws.Range("N58").Value = -1330.5
$ws.Range("H59").Value = 47557.5
$ws.Range("I59").Value = 45000
$ws.Range("J59").Value = 50115
$ws.Range("K59").Value = 45000
$ws.Range("L59").Value = 50115
$ws.Range("M59").Value = -43855
$ws.Range("N59").Value = -52405
$ws.Range("H94").Value = 1470.1428
$ws.Range("J94").Value = 1498.8334
$ws.Range("L94").Value = 1498.8334
$ws.Range("N94").Value = -2400.8334
$ws.Range("H134").Value = 2832.3333
$ws.Range("J134").Value = 2997.5
$ws.Range("L134").Value = 8992.5
$ws.Range("N134").Value = -14062.5
$ws.Range("H136").Value = 946.7778
$ws.Range("I136").Value = 953.1429000000001
$ws.Range("J136").Value = 924.5
$ws.Range("K136").Value = 2859.4287
$ws.Range("L136").Value = 2773.5
$ws.Range("M136").Value = -309.4287000000004
$ws.Range("N136").Value = -7873.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 197.76923
$ws.Range("I17").Value = 32.1
$ws.Range("K17").Value = 96.30000000000001
$ws.Range("M17").Value = 72.69999999999999
$ws.Range("H86").Value = 361.5
$ws.Range("I86").Value = 248
$ws.Range("J86").Value = 475
$ws.Range("K86").Value = 744
$ws.Range("L86").Value = 1425
$ws.Range("M86").Value = 442
$ws.Range("N86").Value = -3797
$ws.Range("H89").Value = 361.5
$ws.Range("I89").Value = 248
$ws.Range("J89").Value = 475
$ws.Range("K89").Value = 2232
$ws.Range("L89").Value = 4275
$ws.Range("M89").Value = 3696
$ws.Range("N89").Value = -16131
$ws.Range("H107").Value = 579.3333
$ws.Range("I107").Value = 377
$ws.Range("K107").Value = 1131
$ws.Range("M107").Value = 789
$ws.Range("H134").Value = 7499.25
$ws.Range("I134").Value = 4999
$ws.Range("K134").Value = 14997
$ws.Range("M134").Value = -9927

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 1231
$ws.Range("I13").Value = 1231
$ws.Range("K13").Value = 1231
$ws.Range("M13").Value = -1092
$ws.Range("H52").Value = 15000
$ws.Range("J52").Value = 15000
$ws.Range("L52").Value = 15000
$ws.Range("N52").Value = -15518

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1390.5714
$ws.Range("J16").Value = 1637.5
$ws.Range("L16").Value = 1637.5
$ws.Range("N16").Value = -1977.5
$ws.Range("H22").Value = 870.7895
$ws.Range("I22").Value = 786.2727
$ws.Range("J22").Value = 987
$ws.Range("K22").Value = 786.2727
$ws.Range("L22").Value = 987
$ws.Range("M22").Value = -491.2727
$ws.Range("N22").Value = -1577
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("H27").Value = 870.7895
$ws.Range("I27").Value = 786.2727
$ws.Range("J27").Value = 987
$ws.Range("K27").Value = 786.2727
$ws.Range("L27").Value = 987
$ws.Range("M27").Value = -679.2727
$ws.Range("N27").Value = -1201
$ws.Range("H132").Value = 3005.25
$ws.Range("J132").Value = 3005.25
$ws.Range("L132").Value = 9015.75
$ws.Range("N132").Value = -14075.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3881.25
$ws.Range("I126").Value = 3077.5
$ws.Range("J126").Value = 7900
$ws.Range("K126").Value = 9232.5
$ws.Range("L126").Value = 23700
$ws.Range("M126").Value = -6762.5
$ws.Range("N126").Value = -28640
$ws.Range("H136").Value = 1532.7693
$ws.Range("I136").Value = 1538.8182
$ws.Range("K136").Value = 4616.4546
$ws.Range("M136").Value = -2066.4546
